# Apply "Optuna Attempt (go back with original)" changes to the forecast workbook.
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---
# Row 2
$wsForecast.Range("D2").Value = 24
$wsForecast.Range("H2").Value = 7.5
$wsForecast.Range("L2").Value = 1.14

# Row 3
$wsForecast.Range("D3").Value = 22
$wsForecast.Range("H3").Value = 7.09
$wsForecast.Range("L3").Value = 1.01

# Row 4
$wsForecast.Range("D4").Value = 21
$wsForecast.Range("H4").Value = 6.38
$wsForecast.Range("L4").Value = 1.19

# Row 5
$wsForecast.Range("D5").Value = 22
$wsForecast.Range("H5").Value = 5.14
$wsForecast.Range("L5").Value = 1.05

# Row 6
$wsForecast.Range("D6").Value = 23
$wsForecast.Range("H6").Value = 3.96
$wsForecast.Range("L6").Value = 1.11

# Row 7
$wsForecast.Range("D7").Value = 19
$wsForecast.Range("H7").Value = 3.58
$wsForecast.Range("J7").Value = "Normal"
$wsForecast.Range("L7").Value = 0.92

# Row 8
$wsForecast.Range("H8").Value = 1.57
$wsForecast.Range("I8").Value = "Low"
$wsForecast.Range("J8").Value = "Normal"
$wsForecast.Range("L8").Value = 1.06

# Row 9
$wsForecast.Range("H9").Value = 0.55
$wsForecast.Range("I9").Value = "Low"
$wsForecast.Range("L9").Value = 0.97

# Row 10
$wsForecast.Range("L10").Value = 0.93

# Row 11
$wsForecast.Range("L11").Value = 1.07

# Row 12
$wsForecast.Range("L12").Value = 1.06

# Row 13
$wsForecast.Range("L13").Value = 1.18

# Row 14
$wsForecast.Range("L14").Value = 0.9

# Row 15
$wsForecast.Range("D15").Value = 22
$wsForecast.Range("L15").Value = 1.06

# Row 16
$wsForecast.Range("L16").Value = 0.87

# Row 17
$wsForecast.Range("L17").Value = 1.1

# --- Summary sheet updates ---
# These cells hold numeric-looking values stored as text, so prefix with an
# apostrophe to force Excel to keep them as text rather than converting to numbers.
$wsSummary.Range("B9").Value = "'425"
$wsSummary.Range("B10").Value = "'195"
$wsSummary.Range("B11").Value = "'89"
$wsSummary.Range("B14").Value = "'19"
